$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 536.8469228506967
$ws.Range("D2").Value = 129.7761258031
$ws.Range("F2").Value = 444
$ws.Range("G2").Value = 490
$ws.Range("H2").Value = 594
$ws.Range("C3").Value = 36.37627465722989
$ws.Range("D3").Value = 6.672834355817785
$ws.Range("E3").Value = 13.96
$ws.Range("F3").Value = 31.37
$ws.Range("G3").Value = 36.49
$ws.Range("H3").Value = 40.78
$ws.Range("C4").Value = 1.922306955007867
$ws.Range("D4").Value = 2.60020301288987
$ws.Range("F4").Value = 0.64
$ws.Range("G4").Value = 1.25
$ws.Range("H4").Value = 2.33
$ws.Range("C5").Value = 322.7852664719634
$ws.Range("D5").Value = 9.604742530486975
$ws.Range("F5").Value = 317.24
$ws.Range("G5").Value = 323.4
$ws.Range("H5").Value = 329.81
$ws.Range("C6").Value = 22.16748994820075
$ws.Range("D6").Value = 2.911911700136196
$ws.Range("F6").Value = 20.31
$ws.Range("G6").Value = 22.08
$ws.Range("H6").Value = 23.97
$ws.Range("I6").Value = 43.95
$ws.Range("C7").Value = -76.7367893705309
$ws.Range("D7").Value = 22.80897747747082
$ws.Range("G7").Value = -74
$ws.Range("C8").Value = 7.415246339073944
$ws.Range("D8").Value = 7.114583036411164
$ws.Range("I8").Value = 25.2
$ws.Range("C9").Value = 9.322031633382492
$ws.Range("D9").Value = 1.685562417434868
$ws.Range("C10").Value = 867.8302641267222
$ws.Range("D10").Value = 0.4616739866011528
$ws.Range("C11").Value = 0.5556674926208094
$ws.Range("D11").Value = 0.5888198148966381
$ws.Range("C12").Value = 22.76645264330656
$ws.Range("D12").Value = 12.29517394858069
$ws.Range("C13").Value = 0.673383608429558
$ws.Range("D13").Value = 0.7508290800637812
$ws.Range("C14").Value = 1.831352490248735
$ws.Range("D14").Value = 1.666409335600012
$ws.Range("C15").Value = 93.99678937053073
$ws.Range("D15").Value = 22.80897747747082
$ws.Range("G15").Value = 91.25999999999999
$ws.Range("C16").Value = -85.86187906500264
$ws.Range("D16").Value = 20.38336693616429
$ws.Range("F16").Value = -102.2324940763249
$ws.Range("H16").Value = -70.26572375596102
$ws.Range("C17").Value = -78.44663272592867
$ws.Range("D17").Value = 25.37424934572375
$ws.Range("F17").Value = -93.41392685158225
$ws.Range("G17").Value = -73.79009749652566
$ws.Range("H17").Value = -60.3175485570292
